$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 67 (weekly price-report update). All rows
# from 67 downward shift to 68 downward, matching the new week's entry
# being prepended to this subset's history.
$ws.Rows.Item(67).Insert()

# Populate the newly inserted row 67 with this week's record (same
# market/product/category dimensions as every other row in this block).
$ws.Cells.Item(67, 1).Value = 11
$ws.Cells.Item(67, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(67, 3).Value = "Bíobío"
$ws.Cells.Item(67, 4).Value = 44460
$ws.Cells.Item(67, 5).Value = 8
$ws.Cells.Item(67, 6).Value = "Fruta"
$ws.Cells.Item(67, 7).Value = 100108
$ws.Cells.Item(67, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(67, 9).Value = 100108005
$ws.Cells.Item(67, 10).Value = "Piña"
$ws.Cells.Item(67, 11).Value = "Caramelo"
$ws.Cells.Item(67, 12).Value = "Segunda"
$ws.Cells.Item(67, 13).Value = 150
$ws.Cells.Item(67, 14).Value = 19000
$ws.Cells.Item(67, 15).Value = 20000
$ws.Cells.Item(67, 16).Value = 19667
$ws.Cells.Item(67, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(67, 18).Value = "Ecuador"
$ws.Cells.Item(67, 19).Value = 1405
$ws.Cells.Item(67, 20).Value = 14
